# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Sat Jun  3 14:22:32 UTC 2023 with GitHub Actions".
# All target cells hold plain text (prices use "." as a thousands-style
# separator rather than being real numbers, and volumes are padded percent
# strings), so each write temporarily forces a Text number format to stop
# Excel from auto-converting the literal into a numeric value, then restores
# the default "Normal" style so the cell format matches the original file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.175.88"
Set-TextValue "E2" "  +0.81%  "
Set-TextValue "D3" "1.901.64"
Set-TextValue "E3" "  +0.92%  "
Set-TextValue "D4" "0.9994"
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "307.37"
Set-TextValue "E5" "  +0.67%  "
Set-TextValue "D6" "0.9995"
Set-TextValue "E6" "  -0.07%  "
Set-TextValue "D7" "0.5196"
Set-TextValue "E7" "  +1.05%  "
Set-TextValue "D8" "0.3767"
Set-TextValue "E8" "  +0.86%  "
Set-TextValue "D9" "0.07271"
Set-TextValue "E9" "  +1.18%  "
Set-TextValue "D10" "21.16"
Set-TextValue "E10" "  +0.50%  "
Set-TextValue "D11" "0.9028"
Set-TextValue "E11" "  +0.51%  "
Set-TextValue "D12" "0.08344"
Set-TextValue "E12" "  +9.14%  "
Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.910.56"
Set-TextValue "E13" "  +1.72%  "
Set-TextValue "B14" "Litecoin"
Set-TextValue "C14" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D14" "96.77"
Set-TextValue "E14" "  +3.50%  "
Set-TextValue "D15" "5.291"
Set-TextValue "E15" "  +1.17%  "
Set-TextValue "D16" "0.9993"
Set-TextValue "E16" "  -0.15%  "
Set-TextValue "D17" "0.000008655"
Set-TextValue "E17" "  +2.05%  "
Set-TextValue "D18" "14.53"
Set-TextValue "E18" "  +1.07%  "
Set-TextValue "D19" "0.9993"
Set-TextValue "E19" "  -0.09%  "
Set-TextValue "D20" "27.217.16"
Set-TextValue "E20" "  +0.77%  "
Set-TextValue "D21" "5.086"
Set-TextValue "E21" "  +0.86%  "
Set-TextValue "D22" "2.134.34"
Set-TextValue "E22" "  +0.12%  "
Set-TextValue "D23" "10.64"
Set-TextValue "E23" "  +1.02%  "
Set-TextValue "D24" "6.439"
Set-TextValue "E24" "  +0.93%  "
Set-TextValue "D25" "2.327"
Set-TextValue "E25" "  +2.16%  "
Set-TextValue "D26" "146.31"
Set-TextValue "E26" "  -0.07%  "
Set-TextValue "D27" "1.751"
Set-TextValue "E27" "  +1.60%  "
Set-TextValue "D28" "18.21"
Set-TextValue "E28" "  +1.10%  "
Set-TextValue "D29" "114.94"
Set-TextValue "E29" "  +0.89%  "
Set-TextValue "D30" "4.822"
Set-TextValue "E30" "  +1.13%  "
Set-TextValue "D31" "4.897"
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "D32" "0.09262"
Set-TextValue "E32" "  +0.97%  "
Set-TextValue "D33" "0.05075"
Set-TextValue "E33" "  +0.90%  "
Set-TextValue "D34" "0.8002"
Set-TextValue "E34" "  +4.77%  "
Set-TextValue "D35" "1.244"
Set-TextValue "E35" "  +1.10%  "
Set-TextValue "D36" "3.417"
Set-TextValue "E36" "  +4.77%  "
Set-TextValue "D37" "2.962"
Set-TextValue "E37" "  -0.43%  "
Set-TextValue "E38" "  +0.45%  "
Set-TextValue "D39" "0.5709"
Set-TextValue "E39" "  +2.24%  "
Set-TextValue "D40" "0.02002"
Set-TextValue "E40" "  +0.81%  "
Set-TextValue "D41" "1.077"
Set-TextValue "E41" "  +0.72%  "
Set-TextValue "D42" "9.037"
Set-TextValue "E42" "  -0.06%  "
Set-TextValue "D43" "6.588"
Set-TextValue "E43" "  -0.29%  "
Set-TextValue "D44" "116.70"
Set-TextValue "D45" "0.1519"
Set-TextValue "E45" "  +1.35%  "
Set-TextValue "D46" "0.4855"
Set-TextValue "E46" "  +1.10%  "
Set-TextValue "B47" "PaxDollar"
Set-TextValue "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D47" "0.9992"
Set-TextValue "E47" "  -0.09%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "10.12"
Set-TextValue "E48" "  -0.45%  "
Set-TextValue "D49" "1.632"
Set-TextValue "E49" "  +2.56%  "
Set-TextValue "D50" "37.77"
Set-TextValue "E50" "  +0.65%  "
Set-TextValue "D51" "63.96"
Set-TextValue "E51" "  +0.31%  "
